$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "主力流入1286亿"
$ws.Range("A3").Value = "主力流出1421亿"
$ws.Range("A4").Value = "主力净流入-134.8亿"
$ws.Range("A5").Value = "超大单418.6亿491.1亿"
$ws.Range("A6").Value = "大单867.8亿930.0亿"
$ws.Range("A7").Value = "中单1259亿1227亿"
$ws.Range("A8").Value = "小单1247亿1144亿"
